$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "From" sender, "Subject", and "Body" cells that hold the
# incoming order email details.
$ws.Range("A2").Value = '"Deepesh Gavali" <deepzgavali143@gmail.com>'
$ws.Range("B2").Value = "Order"
$ws.Range("C2").Value = "Hello`nI would like to order 2 MacBook and 1 iPhone`n"

# The body cell wraps its text onto multiple lines.
$ws.Range("C2").WrapText = $true
$ws.Rows(2).AutoFit()
